$d = $word.ActiveDocument

# The footer block that used to follow "LOB1004: Cálculo II (Requisito fraco)"
# (a blank paragraph, the "Ver no Jupiter ..." line and the "© 2020 ..."
# copyright line) is removed by this site rebuild; the blank paragraph that
# sits right before the trailing page-break paragraph is left untouched.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1004: Cálculo II*") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $blank = $anchor.Next()
    $jupiter = $blank.Next()
    $copyright = $jupiter.Next()

    if (($jupiter.Range.Text -like "*Ver no Jupiter*") -and `
        ($copyright.Range.Text -like "*Contact: luizeleno@usp.br*")) {
        # Delete from the end backwards so the earlier Range objects
        # (whose positions would otherwise shift) stay valid.
        $copyright.Range.Delete()
        $jupiter.Range.Delete()
        $blank.Range.Delete()
    }
}
